# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price / profit recalculation updates
# to the Yojimbo_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 439.41666
$ws.Range("I19").Value = 411.375
$ws.Range("J19").Value = 495.5
$ws.Range("K19").Value = 411.375
$ws.Range("L19").Value = 495.5
$ws.Range("M19").Value = -236.375
$ws.Range("N19").Value = -845.5

$ws.Range("H32").Value = 906.2222
$ws.Range("I32").Value = 850.1667
$ws.Range("J32").Value = 1018.3333
$ws.Range("K32").Value = 850.1667
$ws.Range("L32").Value = 1018.3333
$ws.Range("M32").Value = -524.1667
$ws.Range("N32").Value = -1670.3333

$ws.Range("H43").Value = 5383.6816
$ws.Range("I43").Value = 13276.5
$ws.Range("J43").Value = 873.5
$ws.Range("K43").Value = 13276.5
$ws.Range("L43").Value = 873.5
$ws.Range("M43").Value = -13207.5
$ws.Range("N43").Value = -1011.5

$ws.Range("H51").Value = 1892.2142
$ws.Range("I51").Value = 1869.75
$ws.Range("J51").Value = 1901.2
$ws.Range("K51").Value = 1869.75
$ws.Range("L51").Value = 1901.2
$ws.Range("M51").Value = -1385.75
$ws.Range("N51").Value = -2869.2

$ws.Range("H55").Value = 84.09999999999999
$ws.Range("J55").Value = 90.75
$ws.Range("L55").Value = 90.75
$ws.Range("N55").Value = -518.75

$ws.Range("H98").Value = 1788.0975
$ws.Range("I98").Value = 1201.7241
$ws.Range("J98").Value = 3205.1667
$ws.Range("K98").Value = 1201.7241
$ws.Range("L98").Value = 3205.1667
$ws.Range("M98").Value = 296.2759000000001
$ws.Range("N98").Value = -6201.1667

$ws.Range("H116").Value = 10002661
$ws.Range("I116").Value = 2486.8125
$ws.Range("J116").Value = 27780748
$ws.Range("K116").Value = 2486.8125
$ws.Range("L116").Value = 27780748
$ws.Range("M116").Value = 955.1875
$ws.Range("N116").Value = -27787632

$ws.Range("H122").Value = 1788.0975
$ws.Range("I122").Value = 1201.7241
$ws.Range("J122").Value = 3205.1667
$ws.Range("K122").Value = 3605.1723
$ws.Range("L122").Value = 9615.500100000001
$ws.Range("M122").Value = -1155.1723
$ws.Range("N122").Value = -14515.5001

$ws.Range("H129").Value = 510.66666
$ws.Range("J129").Value = 958.5
$ws.Range("L129").Value = 2875.5
$ws.Range("N129").Value = -12875.5

$ws.Range("H132").Value = 3473600
$ws.Range("I132").Value = 4033425.5
$ws.Range("J132").Value = 2682.4
$ws.Range("K132").Value = 12100276.5
$ws.Range("L132").Value = 8047.200000000001
$ws.Range("M132").Value = -12097746.5
$ws.Range("N132").Value = -13107.2


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 860.7857
$ws.Range("I2").Value = 930.1429000000001
$ws.Range("J2").Value = 791.4286
$ws.Range("K2").Value = 930.1429000000001
$ws.Range("L2").Value = 791.4286
$ws.Range("M2").Value = -817.1429000000001
$ws.Range("N2").Value = -1017.4286

$ws.Range("H32").Value = 3046.3613
$ws.Range("I32").Value = 2193.8684
$ws.Range("J32").Value = 12302
$ws.Range("K32").Value = 2193.8684
$ws.Range("L32").Value = 12302
$ws.Range("M32").Value = -1906.8684
$ws.Range("N32").Value = -12876

$ws.Range("H110").Value = 1036.3077
$ws.Range("I110").Value = 861.0909
$ws.Range("K110").Value = 861.0909
$ws.Range("M110").Value = 1183.9091

$ws.Range("H116").Value = 860.7857
$ws.Range("I116").Value = 930.1429000000001
$ws.Range("J116").Value = 791.4286
$ws.Range("K116").Value = 930.1429000000001
$ws.Range("L116").Value = 791.4286
$ws.Range("M116").Value = 1363.8571
$ws.Range("N116").Value = -5379.4286


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 860.7857
$ws.Range("I3").Value = 930.1429000000001
$ws.Range("J3").Value = 791.4286
$ws.Range("K3").Value = 930.1429000000001
$ws.Range("L3").Value = 791.4286
$ws.Range("M3").Value = -816.1429000000001
$ws.Range("N3").Value = -1019.4286

$ws.Range("H80").Value = 456.08
$ws.Range("I80").Value = 641.375
$ws.Range("J80").Value = 368.88235
$ws.Range("K80").Value = 641.375
$ws.Range("L80").Value = 368.88235
$ws.Range("M80").Value = 356.625
$ws.Range("N80").Value = -2364.88235

$ws.Range("H83").Value = 456.08
$ws.Range("I83").Value = 641.375
$ws.Range("J83").Value = 368.88235
$ws.Range("K83").Value = 3206.875
$ws.Range("L83").Value = 1844.41175
$ws.Range("M83").Value = 1785.125
$ws.Range("N83").Value = -11828.41175

$ws.Range("H107").Value = 1431.8182
$ws.Range("I107").Value = 1435.7142
$ws.Range("J107").Value = 1425
$ws.Range("K107").Value = 1435.7142
$ws.Range("L107").Value = 1425
$ws.Range("M107").Value = 484.2858000000001
$ws.Range("N107").Value = -5265


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1198
$ws.Range("I16").Value = 1312.5
$ws.Range("J16").Value = 740
$ws.Range("K16").Value = 1312.5
$ws.Range("L16").Value = 740
$ws.Range("M16").Value = -1025.5
$ws.Range("N16").Value = -1314

$ws.Range("H86").Value = 2378.2307
$ws.Range("I86").Value = 2040.2609
$ws.Range("J86").Value = 4969.3335
$ws.Range("K86").Value = 2040.2609
$ws.Range("L86").Value = 4969.3335
$ws.Range("M86").Value = -917.2609
$ws.Range("N86").Value = -7215.3335

$ws.Range("H89").Value = 2378.2307
$ws.Range("I89").Value = 2040.2609
$ws.Range("J89").Value = 4969.3335
$ws.Range("K89").Value = 10201.3045
$ws.Range("L89").Value = 24846.6675
$ws.Range("M89").Value = -4585.3045
$ws.Range("N89").Value = -36078.6675

$ws.Range("H99").Value = 2034.7693
$ws.Range("I99").Value = 1956.7778
$ws.Range("J99").Value = 2210.25
$ws.Range("K99").Value = 1956.7778
$ws.Range("L99").Value = 2210.25
$ws.Range("M99").Value = -458.7778000000001
$ws.Range("N99").Value = -5206.25

$ws.Range("H113").Value = 1198
$ws.Range("I113").Value = 1312.5
$ws.Range("J113").Value = 740
$ws.Range("K113").Value = 1312.5
$ws.Range("L113").Value = 740
$ws.Range("M113").Value = 857.5
$ws.Range("N113").Value = -5080

$ws.Range("H122").Value = 2847.9412
$ws.Range("I122").Value = 2969.6875
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 8909.0625
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -6459.0625
$ws.Range("N122").Value = -7600

$ws.Range("H126").Value = 2034.7693
$ws.Range("I126").Value = 1956.7778
$ws.Range("J126").Value = 2210.25
$ws.Range("K126").Value = 5870.3334
$ws.Range("L126").Value = 6630.75
$ws.Range("M126").Value = -3400.3334
$ws.Range("N126").Value = -11570.75

$ws.Range("H134").Value = 1373.75
$ws.Range("I134").Value = 1285.4584
$ws.Range("J134").Value = 1903.5
$ws.Range("K134").Value = 3856.3752
$ws.Range("L134").Value = 5710.5
$ws.Range("M134").Value = -1321.3752
$ws.Range("N134").Value = -10780.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 491
$ws.Range("I2").Value = 38.75
$ws.Range("J2").Value = 749.4286
$ws.Range("K2").Value = 232.5
$ws.Range("L2").Value = 4496.571599999999
$ws.Range("M2").Value = -119.5
$ws.Range("N2").Value = -4722.571599999999

$ws.Range("H23").Value = 76.666664
$ws.Range("J23").Value = 95.77778000000001
$ws.Range("L23").Value = 287.33334
$ws.Range("N23").Value = -757.33334

$ws.Range("H68").Value = 10880.3
$ws.Range("J68").Value = 15371.857
$ws.Range("L68").Value = 46115.571
$ws.Range("N68").Value = -47737.571

$ws.Range("H71").Value = 10880.3
$ws.Range("J71").Value = 15371.857
$ws.Range("L71").Value = 138346.713
$ws.Range("N71").Value = -146458.713

$ws.Range("H80").Value = 4595.6
$ws.Range("I80").Value = 4326
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 12978
$ws.Range("L80").Value = 15000
$ws.Range("M80").Value = -12042
$ws.Range("N80").Value = -16872

$ws.Range("H83").Value = 4595.6
$ws.Range("I83").Value = 4326
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 38934
$ws.Range("L83").Value = 45000
$ws.Range("M83").Value = -34254
$ws.Range("N83").Value = -54360

$ws.Range("H92").Value = 425.25
$ws.Range("I92").Value = 299.66666
$ws.Range("J92").Value = 500.6
$ws.Range("K92").Value = 898.9999799999999
$ws.Range("L92").Value = 1501.8
$ws.Range("M92").Value = 349.0000200000001
$ws.Range("N92").Value = -3997.8

$ws.Range("H98").Value = 667487.0600000001
$ws.Range("I98").Value = 620.6
$ws.Range("J98").Value = 1000920.3
$ws.Range("K98").Value = 1861.8
$ws.Range("L98").Value = 3002760.9
$ws.Range("M98").Value = -363.8000000000002
$ws.Range("N98").Value = -3005756.9

$ws.Range("H122").Value = 1588.4
$ws.Range("I122").Value = 1619.8
$ws.Range("J122").Value = 1572.7
$ws.Range("K122").Value = 14578.2
$ws.Range("L122").Value = 14154.3
$ws.Range("M122").Value = -12128.2
$ws.Range("N122").Value = -19054.3


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4751.839
$ws.Range("I70").Value = 4679.316
$ws.Range("J70").Value = 4866.6665
$ws.Range("K70").Value = 4679.316
$ws.Range("L70").Value = 4866.6665
$ws.Range("M70").Value = -4409.316
$ws.Range("N70").Value = -5406.6665

$ws.Range("H73").Value = 4751.839
$ws.Range("I73").Value = 4679.316
$ws.Range("J73").Value = 4866.6665
$ws.Range("K73").Value = 4679.316
$ws.Range("L73").Value = 4866.6665
$ws.Range("M73").Value = -3743.316
$ws.Range("N73").Value = -6738.6665

$ws.Range("H113").Value = 1012
$ws.Range("I113").Value = 1011
$ws.Range("J113").Value = 1013
$ws.Range("K113").Value = 1011
$ws.Range("L113").Value = 1013
$ws.Range("M113").Value = 1159
$ws.Range("N113").Value = -5353

$ws.Range("H122").Value = 2941.9565
$ws.Range("I122").Value = 2740.5557
$ws.Range("J122").Value = 3071.4285
$ws.Range("K122").Value = 8221.667099999999
$ws.Range("L122").Value = 9214.2855
$ws.Range("M122").Value = -5771.667099999999
$ws.Range("N122").Value = -14114.2855

$ws.Range("H126").Value = 51136.45
$ws.Range("I126").Value = 77963.766
$ws.Range("J126").Value = 1314.2858
$ws.Range("K126").Value = 233891.298
$ws.Range("L126").Value = 3942.8574
$ws.Range("M126").Value = -231421.298
$ws.Range("N126").Value = -8882.857400000001

$ws.Range("H132").Value = 2626.4119
$ws.Range("I132").Value = 2395.077
$ws.Range("J132").Value = 3378.25
$ws.Range("K132").Value = 7185.231000000001
$ws.Range("L132").Value = 10134.75
$ws.Range("M132").Value = -4655.231000000001
$ws.Range("N132").Value = -15194.75


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 384.14285
$ws.Range("I22").Value = 384.14285
$ws.Range("K22").Value = 384.14285
$ws.Range("M22").Value = -89.14285000000001

$ws.Range("H27").Value = 384.14285
$ws.Range("I27").Value = 384.14285
$ws.Range("K27").Value = 384.14285
$ws.Range("M27").Value = -277.14285

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""

$ws.Range("H122").Value = 3530.7144
$ws.Range("I122").Value = 3612.2727
$ws.Range("J122").Value = 3441
$ws.Range("K122").Value = 10836.8181
$ws.Range("L122").Value = 10323
$ws.Range("M122").Value = -8386.8181
$ws.Range("N122").Value = -15223

$ws.Range("H136").Value = 2685
$ws.Range("I136").Value = 1872.7451
$ws.Range("J136").Value = 4657.619
$ws.Range("K136").Value = 5618.2353
$ws.Range("L136").Value = 13972.857
$ws.Range("M136").Value = -3068.2353
$ws.Range("N136").Value = -19072.857

